# GamingWebsite.pptx edit script
# - Slide 1 title: "Buddies Hubs" -> "Buddies Hub"
# - Handout master date placeholder: 8/27/2021 -> 10/1/2021
# - Notes master date placeholder: 8/27/2021 -> 10/1/2021

$p = $ppt.ActivePresentation

# --- 1. Fix the title text on the first slide ("content removed": the
#        trailing "s" in "Buddies Hubs" is removed) -----------------------
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item("Title 2")
$titleShape.TextFrame.TextRange.Text = "Buddies Hub"

# --- 2. Refresh the auto date placeholders (handout master + notes master)
#        so the cached "datetimeFigureOut" text matches the later save date
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "10/1/2021"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "10/1/2021"
